$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.027520695437743
$ws.Range("D2").Value = 1.032301774436475
$ws.Range("E2").Value = 1.027609874303749
$ws.Range("I2").Value = 1.033536576406378
$ws.Range("J2").Value = 1.032678109988713
$ws.Range("K2").Value = 1.035107398427057
$ws.Range("L2").Value = 1.030429102474869
$ws.Range("N2").Value = 1.034144631354739

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.028301072450298
$ws.Range("D3").Value = 1.032878511955579
$ws.Range("E3").Value = 1.02826801873824
$ws.Range("I3").Value = 1.03368227300884
$ws.Range("J3").Value = 1.033099575735387
$ws.Range("K3").Value = 1.035493500853492
$ws.Range("L3").Value = 1.030895412652906
$ws.Range("N3").Value = 1.034566695631109

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.028806608519734
$ws.Range("D4").Value = 1.033252132454826
$ws.Range("E4").Value = 1.028694776556743
$ws.Range("I4").Value = 1.033775613710888
$ws.Range("J4").Value = 1.033372196915805
$ws.Range("K4").Value = 1.035743066147848
$ws.Range("L4").Value = 1.031197356885715
$ws.Range("N4").Value = 1.034839703964874

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.029019272836835
$ws.Range("D5").Value = 1.033409304191379
$ws.Range("E5").Value = 1.028874397673045
$ws.Range("I5").Value = 1.033814629610479
$ws.Range("J5").Value = 1.033486782661957
$ws.Range("K5").Value = 1.035847917666564
$ws.Range("L5").Value = 1.031324343115634
$ws.Range("N5").Value = 1.034954452435928

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.029054988087535
$ws.Range("D6").Value = 1.033435699923981
$ws.Range("E6").Value = 1.028904569230505
$ws.Range("I6").Value = 1.033821167363204
$ws.Range("J6").Value = 1.033506020653993
$ws.Range("K6").Value = 1.035865518803235
$ws.Range("L6").Value = 1.03134566748513
$ws.Range("N6").Value = 1.03497371774812

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.028809449614982
$ws.Range("D7").Value = 1.03325423219416
$ws.Range("E7").Value = 1.028697175831864
$ws.Range("I7").Value = 1.033776135926442
$ws.Range("J7").Value = 1.033373728112961
$ws.Range("K7").Value = 1.035744467438987
$ws.Range("L7").Value = 1.031199053491763
$ws.Range("N7").Value = 1.034841237336505

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.027784306281523
$ws.Range("D8").Value = 1.032496594819529
$ws.Range("E8").Value = 1.027832110922644
$ws.Range("I8").Value = 1.033586008347227
$ws.Range("J8").Value = 1.032820565226457
$ws.Range("K8").Value = 1.035237938375111
$ws.Range("L8").Value = 1.030586649731124
$ws.Range("N8").Value = 1.034287288895253

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.02598239819471
$ws.Range("D9").Value = 1.031164940409345
$ws.Range("E9").Value = 1.026314693355839
$ws.Range("I9").Value = 1.033243853563985
$ws.Range("J9").Value = 1.031845152066742
$ws.Range("K9").Value = 1.034343371866929
$ws.Range("L9").Value = 1.029509192962856
$ws.Range("N9").Value = 1.03331049053697

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.024784281678368
$ws.Range("D10").Value = 1.030279571377919
$ws.Range("E10").Value = 1.025307864385804
$ws.Range("I10").Value = 1.03301100448552
$ws.Range("J10").Value = 1.031194508535528
$ws.Range("K10").Value = 1.033745733651523
$ws.Range("L10").Value = 1.028792104773625
$ws.Range("N10").Value = 1.032658923017312

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.02426625629888
$ws.Range("D11").Value = 1.02989679057021
$ws.Range("E11").Value = 1.024873055851738
$ws.Range("I11").Value = 1.032909062120811
$ws.Range("J11").Value = 1.030912702751864
$ws.Range("K11").Value = 1.033486668187737
$ws.Range("L11").Value = 1.028481904987323
$ws.Range("N11").Value = 1.032376717037112

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.024073955764672
$ws.Range("D12").Value = 1.029754699333201
$ws.Range("E12").Value = 1.024711724186278
$ws.Range("I12").Value = 1.032871029101589
$ws.Range("J12").Value = 1.030808017995343
$ws.Range("K12").Value = 1.033390398456624
$ws.Range("L12").Value = 1.028366730194739
$ws.Range("N12").Value = 1.032271883616228

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.024115199549496
$ws.Range("D13").Value = 1.029785174254312
$ws.Range("E13").Value = 1.024746322397756
$ws.Range("I13").Value = 1.032879194862144
$ws.Range("J13").Value = 1.030830473640227
$ws.Range("K13").Value = 1.033411050486251
$ws.Range("L13").Value = 1.02839143341743
$ws.Range("N13").Value = 1.032294371150704

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.024250358275622
$ws.Range("D14").Value = 1.029885043404599
$ws.Range("E14").Value = 1.024859716532243
$ws.Range("I14").Value = 1.032905921703698
$ws.Range("J14").Value = 1.030904049663762
$ws.Range("K14").Value = 1.033478711335739
$ws.Range("L14").Value = 1.028472383630948
$ws.Range("N14").Value = 1.032368051660632

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.024333649609663
$ws.Range("D15").Value = 1.029946588148598
$ws.Range("E15").Value = 1.024929605728804
$ws.Range("I15").Value = 1.032922366875029
$ws.Range("J15").Value = 1.030949381057819
$ws.Range("K15").Value = 1.033520393989972
$ws.Range("L15").Value = 1.028522266031901
$ws.Range("N15").Value = 1.032413447430468

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.024818677640542
$ws.Range("D16").Value = 1.030304987875841
$ws.Range("E16").Value = 1.025336745706958
$ws.Range("I16").Value = 1.033017746591506
$ws.Range("J16").Value = 1.031213209642353
$ws.Range("K16").Value = 1.033762921106411
$ws.Range("L16").Value = 1.028812698249956
$ws.Range("N16").Value = 1.032677650681854

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.025123129611886
$ws.Range("D17").Value = 1.030529961870703
$ws.Range("E17").Value = 1.02559244454807
$ws.Range("I17").Value = 1.033077277201979
$ws.Range("J17").Value = 1.031378683832033
$ws.Range("K17").Value = 1.033914976874395
$ws.Range("L17").Value = 1.028994961210851
$ws.Range("N17").Value = 1.032843359863863

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.025300785168862
$ws.Range("D18").Value = 1.030661242172459
$ws.Range("E18").Value = 1.02574170066471
$ws.Range("I18").Value = 1.033111892557661
$ws.Range("J18").Value = 1.031475195005187
$ws.Range("K18").Value = 1.034003640918216
$ws.Range("L18").Value = 1.029101301354294
$ws.Range("N18").Value = 1.032940008093953

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.025361373594297
$ws.Range("D19").Value = 1.030706014932689
$ws.Range("E19").Value = 1.025792611960816
$ws.Range("I19").Value = 1.03312367718447
$ws.Range("J19").Value = 1.031508101570444
$ws.Range("K19").Value = 1.034033868356237
$ws.Range("L19").Value = 1.029137565508582
$ws.Range("N19").Value = 1.032972961390307

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.025090457141737
$ws.Range("D20").Value = 1.030505818403187
$ws.Range("E20").Value = 1.025564998954362
$ws.Range("I20").Value = 1.033070901280038
$ws.Range("J20").Value = 1.031360930744068
$ws.Range("K20").Value = 1.033898665565162
$ws.Range("L20").Value = 1.028975403084414
$ws.Range("N20").Value = 1.032825581564477

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.024210554149939
$ws.Range("D21").Value = 1.029855631912149
$ws.Range("E21").Value = 1.024826319922492
$ws.Range("I21").Value = 1.032898055921315
$ws.Range("J21").Value = 1.030882383610789
$ws.Range("K21").Value = 1.033458788030259
$ws.Range("L21").Value = 1.028448544494949
$ws.Range("N21").Value = 1.032346354839379

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.023658003003678
$ws.Range("D22").Value = 1.029447359011825
$ws.Range("E22").Value = 1.024362899983366
$ws.Range("I22").Value = 1.032788414851583
$ws.Range("J22").Value = 1.030581447526175
$ws.Range("K22").Value = 1.033181981411915
$ws.Range("L22").Value = 1.028117561541691
$ws.Range("N22").Value = 1.032044991391007

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.023950855835906
$ws.Range("D23").Value = 1.029663741761482
$ws.Range("E23").Value = 1.024608470588753
$ws.Range("I23").Value = 1.032846629001183
$ws.Range("J23").Value = 1.030740984137351
$ws.Range("K23").Value = 1.033328743941221
$ws.Range("L23").Value = 1.028292995390668
$ws.Range("N23").Value = 1.03220475456247

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025105220196619
$ws.Range("D24").Value = 1.030516727621749
$ws.Range("E24").Value = 1.025577400093128
$ws.Range("I24").Value = 1.033073782618027
$ws.Range("J24").Value = 1.031368952622536
$ws.Range("K24").Value = 1.033906036028888
$ws.Range("L24").Value = 1.02898424046973
$ws.Range("N24").Value = 1.032833614834932

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.026447686710587
$ws.Range("D25").Value = 1.031508790279599
$ws.Range("E25").Value = 1.026706147838568
$ws.Range("I25").Value = 1.033333148572686
$ws.Range("J25").Value = 1.032097390628995
$ws.Range("K25").Value = 1.034574866967415
$ws.Range("L25").Value = 1.029787533477497
$ws.Range("N25").Value = 1.033563087306913
